$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.002.31"
Set-TextValue "E2" "  -0.18%  "
Set-TextValue "D3" "1.631.03"
Set-TextValue "E3" "  -0.87%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.14%  "
Set-TextValue "D5" "214.45"
Set-TextValue "E5" "  -0.77%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.17%  "
Set-TextValue "E8" "  -1.96%  "
Set-TextValue "D9" "0.0620"
Set-TextValue "E9" "  -3.14%  "
Set-TextValue "D10" "18.53"
Set-TextValue "E10" "  -5.49%  "
Set-TextValue "E11" "  -0.96%  "
Set-TextValue "D12" "1.854.44"
Set-TextValue "E12" "  -0.97%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.19"
Set-TextValue "E13" "  -2.11%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.619.87"
Set-TextValue "E14" "  -2.21%  "
Set-TextValue "D15" "0.530"
Set-TextValue "E15" "  -2.80%  "
Set-TextValue "D16" "26.004.35"
Set-TextValue "E16" "  -0.20%  "
Set-TextValue "D17" "0.0₃0741"
Set-TextValue "E17" "  -2.97%  "
Set-TextValue "D18" "61.53"
Set-TextValue "E18" "  -3.08%  "
Set-TextValue "D19" "1.00"
Set-TextValue "E19" "  -0.14%  "
Set-TextValue "D20" "193.86"
Set-TextValue "E20" "  -0.36%  "
Set-TextValue "D21" "4.26"
Set-TextValue "E21" "  -2.44%  "
Set-TextValue "E22" "  -3.77%  "
Set-TextValue "D23" "6.09"
Set-TextValue "E23" "  -2.04%  "
Set-TextValue "E24" "  +1.03%  "
Set-TextValue "D25" "144.12"
Set-TextValue "E25" "  +0.54%  "
Set-TextValue "E26" "  -0.09%  "
Set-TextValue "D27" "1.72"
Set-TextValue "E27" "  -4.43%  "
Set-TextValue "D28" "6.75"
Set-TextValue "E28" "  -2.00%  "
Set-TextValue "D29" "15.28"
Set-TextValue "E29" "  -1.61%  "
Set-TextValue "E30" "  -1.29%  "
Set-TextValue "E31" "  -2.24%  "
Set-TextValue "D32" "3.14"
Set-TextValue "E32" "  -3.95%  "
Set-TextValue "D33" "3.12"
Set-TextValue "E33" "  -5.39%  "
Set-TextValue "E34" "  -2.78%  "
Set-TextValue "E35" "  -2.70%  "
Set-TextValue "D36" "1.125.18"
Set-TextValue "E36" "  -0.60%  "
Set-TextValue "D37" "0.853"
Set-TextValue "E37" "  -5.93%  "
Set-TextValue "E39" "  -3.23%  "
Set-TextValue "E40" "  -2.18%  "
Set-TextValue "D41" "98.23"
Set-TextValue "E41" "  -0.85%  "
Set-TextValue "E42" "  -3.69%  "
Set-TextValue "D43" "1.764.93"
Set-TextValue "E43" "  -0.94%  "
Set-TextValue "E44" "  -5.44%  "
Set-TextValue "D45" "0.0₆0111"
Set-TextValue "E45" "  -4.80%  "
Set-TextValue "E46" "  +2.01%  "
Set-TextValue "D47" "54.57"
Set-TextValue "E47" "  -3.58%  "
Set-TextValue "E48" "  -0.70%  "
Set-TextValue "D49" "0.412"
Set-TextValue "E49" "  -0.62%  "
Set-TextValue "D50" "7.48"
Set-TextValue "E50" "  -4.10%  "
Set-TextValue "D51" "1.00"
Set-TextValue "E51" "  -0.02%  "
